$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for rule R10 (cell E8) to reflect the git update
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the selection state recorded in the saved file (active cell E8)
$ws.Range("E8").Select()
